# Two groups (5B2EP and 3B1EP) of scenarios are run.
# The original table had 5 scenario rows (A2:A6 = 1..5); now only scenarios
# 2, 4 and 5 remain, so the rows that held the old scenarios 1 and 3 are
# removed entirely (the rows below shift up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the worksheet rows that held the old scenario "1" and "3" records
# (sheet rows 4 then 2, bottom-to-top so row numbers stay valid).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

# Re-apply the AutoFilter so its range shrinks to the new data extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:P4").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=OperationScenario_Component_Bui!`$A`$1:`$P`$4"
    }
}

# Restore the selection that was active after the edit.
$ws.Range("C9").Select()
